$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value  = 2.496985546467907
$ws.Range("D6").Value  = 0.07635665951707242
$ws.Range("D7").Value  = -0.002535233822187239
$ws.Range("D8").Value  = 0.5400231511551582
$ws.Range("D9").Value  = 2.533663549832727
$ws.Range("D10").Value = 0.3030152403641788
$ws.Range("D11").Value = 2.417443432035182
$ws.Range("D12").Value = 0.006475575292692774
$ws.Range("D13").Value = 0.2025111021396851
$ws.Range("D14").Value = 0.3792894630601535
$ws.Range("D15").Value = 0.3917238595074686
$ws.Range("D16").Value = 0.03537009864606364
$ws.Range("D17").Value = 0.01064639444370955
$ws.Range("D18").Value = -0.06424314580939124
$ws.Range("D19").Value = -0.02792184497684371
$ws.Range("D21").Value = 0
$ws.Range("D22").Value = 0.5067639732573387
$ws.Range("D23").Value = 0.2651495670929225
